$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new shared string "Сумм" used for header L2, mirroring K2 formatting
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("L2").Value = "Сумм"

# 2. Zero-out the cells that previously held a "2" (homework/lab grade removed)
$ws.Range("G4:J4").Value = 0
$ws.Range("H6:J6").Value = 0
$ws.Range("C7:J7").Value = 0
$ws.Range("G8:J8").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("I11:J11").Value = 0
$ws.Range("I14:J14").Value = 0
$ws.Range("I15:J15").Value = 0
$ws.Range("C17:J17").Value = 0
$ws.Range("G18:J18").Value = 0
$ws.Range("G19:J19").Value = 0
$ws.Range("C21:J21").Value = 0
$ws.Range("E23:J23").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("C26:J26").Value = 0
$ws.Range("C28:J28").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("C30:J30").Value = 0
$ws.Range("F31:J31").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 0

# 3. Populate the new "M" (attendance/М) column, cloning the J-column cell style (s=7)
$ws.Range("J7").Copy()
$ws.Range("M4:M32").PasteSpecial(-4122)
$ws.Range("M4").Value = 4
$ws.Range("M5").Value = 5
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 3
$ws.Range("M8").Value = 4
$ws.Range("M9").Value = 5
$ws.Range("M10").Value = 5
$ws.Range("M11").Value = 5
$ws.Range("M12").Value = 5
$ws.Range("M13").Value = 5
$ws.Range("M14").Value = 5
$ws.Range("M15").Value = 5
$ws.Range("M16").Value = 5
$ws.Range("M17").Value = 3
$ws.Range("M18").Value = 4
$ws.Range("M19").Value = 4
$ws.Range("M20").Value = 5
$ws.Range("M21").Value = 3
$ws.Range("M22").Value = 5
$ws.Range("M23").Value = 3
$ws.Range("M24").Value = 5
$ws.Range("M25").Value = 5
$ws.Range("M26").Value = 3
$ws.Range("M27").Value = 5
$ws.Range("M28").Value = 3
$ws.Range("M29").Value = 5
$ws.Range("M30").Value = 3
$ws.Range("M31").Value = 3
$ws.Range("M32").Value = 5

# 4. Resize columns C:J and L to the new narrower width
$ws.Range("C1:J1").EntireColumn.ColumnWidth = 4.0
$ws.Range("L1").EntireColumn.ColumnWidth = 4.0

# 5. Rebuild conditional formatting: drop old rules, add the 3 new color-scale rules
$ws.Range("J4:J32").FormatConditions.Delete()
$ws.Range("L4:L32").FormatConditions.Delete()
$ws.Range("E25").FormatConditions.Delete()
$ws.Range("G25").FormatConditions.Delete()

$fcC = $ws.Range("C4:J32").FormatConditions.AddColorScale(3)
$fcL = $ws.Range("L4:L32").FormatConditions.AddColorScale(3)
$fcM = $ws.Range("M4:M32").FormatConditions.AddColorScale(3)
$fcC.Priority = 3
$fcL.Priority = 2
$fcM.Priority = 1

# 6. Update the active selection / scroll position to match the edited workbook
$ws.Range("M11").Select() | Out-Null
